$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Autofit the columns that end up with explicit widths in the final sheet
$ws.Columns("A:A").AutoFit()
$ws.Columns("D:D").AutoFit()
$ws.Columns("G:G").AutoFit()

# Row 12
$ws.Range("A12").Value = "b9407f30-f5f8-466e-aff9-25556b57fe6d"
$ws.Range("B12").Value = 10345
$ws.Range("C12").Value = 19843
$ws.Range("D12").Value = "Welcome to St.Loius Airport "
$ws.Range("E12").Value = 123.5
$ws.Range("F12").Value = 206
$ws.Range("G12").Value = "St.Loius Airport shuttle 1"

# Row 13
$ws.Range("A13").Value = "50765cb7-d9ea-4e21-99a4-fa879613a492"
$ws.Range("B13").Value = 62477
$ws.Range("C13").Value = 47058
$ws.Range("D13").Value = "Welcome to St.Loius Airport gateway 2"
$ws.Range("E13").Value = 109.5
$ws.Range("F13").Value = 200
$ws.Range("G13").Value = "St.Loius Airport shuttle 2"

# Apply left-alignment style to the two new rows (creates new cellXfs entry)
$ws.Range("A12:G13").HorizontalAlignment = -4131

# Update selection to mimic post-edit cursor location
$ws.Range("G23").Select()
